# Updated symbol list (applies data refresh captured in the commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "272.97"
Set-TextValue "E2" "0.95%"
Set-TextValue "D3" "26.83"
Set-TextValue "E3" "0.38%"
Set-TextValue "D4" "4.900"
Set-TextValue "E4" "3.86%"
Set-TextValue "D5" "0.06318"
Set-TextValue "E5" "3.13%"
Set-TextValue "D6" "6.906"
Set-TextValue "E6" "2.41%"
Set-TextValue "D7" "3.353"
Set-TextValue "E7" "5.41%"
Set-TextValue "D8" "1.410"
Set-TextValue "E8" "56.91%"
Set-TextValue "D9" "0.8830"
Set-TextValue "E9" "3.21%"
Set-TextValue "D10" "0.1461"
Set-TextValue "E10" "2.26%"
Set-TextValue "D11" "0.05102"
Set-TextValue "E11" "0.64%"
Set-TextValue "D12" "0.07404"
Set-TextValue "E12" "3.48%"
Set-TextValue "D13" "0.03157"
Set-TextValue "E13" "-1.02%"
Set-TextValue "D14" "0.09037"
Set-TextValue "E14" "-0.07%"
Set-TextValue "D15" "0.001568"
Set-TextValue "E15" "2.61%"
Set-TextValue "D16" "0.0006299"
Set-TextValue "E16" "3.73%"
Set-TextValue "D17" "0.006016"
Set-TextValue "E17" "-0.17%"
Set-TextValue "D18" "3.468"
Set-TextValue "E18" "0.14%"
Set-TextValue "D19" "2.272"
Set-TextValue "E19" "-0.28%"
Set-TextValue "E20" "2.51%"
Set-TextValue "E21" "4.06%"
Set-TextValue "D22" "3.905"
Set-TextValue "E22" "1.75%"
Set-TextValue "D23" "0.04342"
Set-TextValue "E23" "1.91%"
Set-TextValue "D24" "0.001176"
Set-TextValue "E24" "-0.28%"
Set-TextValue "D25" "0.003653"
Set-TextValue "E25" "-11.94%"
Set-TextValue "D26" "0.0001199"
Set-TextValue "E26" "-0.20%"
Set-TextValue "E27" "1.12%"
Set-TextValue "D40" "0.04043"
Set-TextValue "E40" "1.84%"
Set-TextValue "D41" "0.006602"
Set-TextValue "E41" "57.29%"
Set-TextValue "D42" "0.1163"
Set-TextValue "E42" "3.85%"
Set-TextValue "D43" "0.002128"
Set-TextValue "E43" "4.34%"
Set-TextValue "D44" "0.01256"
Set-TextValue "E44" "6.60%"
Set-TextValue "D45" "0.00005341"
Set-TextValue "E45" "3.91%"
Set-TextValue "E46" "159.65%"
Set-TextValue "D47" "0.02120"
Set-TextValue "E47" "-29.21%"

$wb.Save()
